$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Тест по выборке из 5000")

$ws.Range("A6").Value = "SVM (SVC) (лидер)"
$ws.Range("B6").Value = "Rbf, C=1, gamma=0.0001"
$ws.Range("C6").Value = 5000
$ws.Range("D6").Value = "CV, 5"
$ws.Range("E6").Value = 0.78
$ws.Range("F6").Value = 0.02
$ws.Range("G6").Value = 10.199999999999999

$ws.Range("A5:G5").Copy() | Out-Null
$ws.Range("A6:G6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("G7").Select() | Out-Null
